$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = -0.4828312689463684
$ws.Range("J3").Value = 0.2158793921891626
$ws.Range("K3").Value = -0.519189168748834
$ws.Range("L3").Value = 2.856722983783609
